# Commit: "data up otto 14th" -- extend community-state survey data
# through 14 Sep 2020 and correct a handful of previously-entered values
# for 03-07 Sep 2020 (rows 217-220).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small numeric corrections to already-present rows 217-220 ---
$ws.Range("L217").Value2 = 19.8429203
$ws.Range("X217").Value2 = 12.9750876
$ws.Range("BC217").Value2 = 19.4970778
$ws.Range("L218").Value2 = 19.6319525
$ws.Range("X218").Value2 = 12.8033247
$ws.Range("BC218").Value2 = 19.7464277
$ws.Range("L219").Value2 = 19.6211898
$ws.Range("X219").Value2 = 13.2547118
$ws.Range("BC219").Value2 = 19.8589155
$ws.Range("G220").Value2 = 15.1351926
$ws.Range("L220").Value2 = 19.4144043
$ws.Range("X220").Value2 = 13.270305
$ws.Range("AW220").Value2 = 21.723013
$ws.Range("BC220").Value2 = 19.8431852

# --- New daily rows: 07 Sep 2020 (data for previously date-only row 221) ---
# --- through 14 Sep 2020 (row 228, date only -- no survey responses yet) ---
# Row 221
$ws.Range("B221").Value2 = 19.7522816
$ws.Range("C221").Value2 = 27.9346891
$ws.Range("D221").Value2 = 26.7135996
$ws.Range("F221").Value2 = 16.5240542
$ws.Range("G221").Value2 = 14.8505697
$ws.Range("H221").Value2 = 15.8369119
$ws.Range("I221").Value2 = 9.507913
$ws.Range("J221").Value2 = 13.2520325
$ws.Range("K221").Value2 = 12.2769754
$ws.Range("L221").Value2 = 19.3741098
$ws.Range("M221").Value2 = 26.147733
$ws.Range("O221").Value2 = 16.0685805
$ws.Range("P221").Value2 = 29.9192125
$ws.Range("Q221").Value2 = 26.718357
$ws.Range("R221").Value2 = 20.4760189
$ws.Range("S221").Value2 = 24.8268822
$ws.Range("T221").Value2 = 26.2297664
$ws.Range("U221").Value2 = 23.4796748
$ws.Range("V221").Value2 = 23.8256672
$ws.Range("W221").Value2 = 10.7808555
$ws.Range("X221").Value2 = 13.2563923
$ws.Range("Y221").Value2 = 10.4354753
$ws.Range("Z221").Value2 = 14.8972611
$ws.Range("AA221").Value2 = 17.5653183
$ws.Range("AB221").Value2 = 28.5258923
$ws.Range("AD221").Value2 = 29.0379352
$ws.Range("AE221").Value2 = 22.1402825
$ws.Range("AF221").Value2 = 18.1970603
$ws.Range("AG221").Value2 = 30.2763785
$ws.Range("AH221").Value2 = 26.9418788
$ws.Range("AI221").Value2 = 10.0938967
$ws.Range("AJ221").Value2 = 10.1021682
$ws.Range("AK221").Value2 = 15.2133472
$ws.Range("AL221").Value2 = 19.2564634
$ws.Range("AM221").Value2 = 10.4623183
$ws.Range("AN221").Value2 = 17.7094465
$ws.Range("AO221").Value2 = 27.8464286
$ws.Range("AP221").Value2 = 12.8057871
$ws.Range("AQ221").Value2 = 13.4819896
$ws.Range("AS221").Value2 = 10.8791124
$ws.Range("AT221").Value2 = 21.6461165
$ws.Range("AU221").Value2 = 33.238607
$ws.Range("AV221").Value2 = 26.4411334
$ws.Range("AW221").Value2 = 21.397675
$ws.Range("AX221").Value2 = 23.2004779
$ws.Range("AY221").Value2 = 17.5284504
$ws.Range("BA221").Value2 = 8.090029599999999
$ws.Range("BB221").Value2 = 14.9033254
$ws.Range("BC221").Value2 = 19.8947684
$ws.Range("BD221").Value2 = 20.2280164
$ws.Range("BE221").Value2 = 20.0201613

# Row 222
$ws.Range("A222").Value2 = "08 09 2020"
$ws.Range("B222").Value2 = 19.5364238
$ws.Range("C222").Value2 = 27.1374849
$ws.Range("D222").Value2 = 26.4325128
$ws.Range("F222").Value2 = 16.2977241
$ws.Range("G222").Value2 = 14.6924722
$ws.Range("H222").Value2 = 16.1604681
$ws.Range("I222").Value2 = 9.4618494
$ws.Range("J222").Value2 = 12.7467105
$ws.Range("K222").Value2 = 12.2769754
$ws.Range("L222").Value2 = 19.3051469
$ws.Range("M222").Value2 = 26.0909485
$ws.Range("O222").Value2 = 15.5510204
$ws.Range("P222").Value2 = 29.7616337
$ws.Range("Q222").Value2 = 26.7433639
$ws.Range("R222").Value2 = 20.5054219
$ws.Range("S222").Value2 = 24.9490745
$ws.Range("T222").Value2 = 26.2575676
$ws.Range("U222").Value2 = 23.1158686
$ws.Range("V222").Value2 = 23.0588734
$ws.Range("W222").Value2 = 10.5335076
$ws.Range("X222").Value2 = 13.3665715
$ws.Range("Y222").Value2 = 10.9951846
$ws.Range("Z222").Value2 = 14.8606901
$ws.Range("AA222").Value2 = 17.5451961
$ws.Range("AB222").Value2 = 28.6765605
$ws.Range("AD222").Value2 = 28.9517722
$ws.Range("AE222").Value2 = 22.4382721
$ws.Range("AF222").Value2 = 17.9522769
$ws.Range("AG222").Value2 = 31.226091
$ws.Range("AH222").Value2 = 26.1622499
$ws.Range("AI222").Value2 = 9.0207557
$ws.Range("AJ222").Value2 = 10.0462817
$ws.Range("AK222").Value2 = 15.7632134
$ws.Range("AL222").Value2 = 19.3904884
$ws.Range("AM222").Value2 = 10.3709032
$ws.Range("AN222").Value2 = 17.6245033
$ws.Range("AO222").Value2 = 27.9413851
$ws.Range("AP222").Value2 = 12.6790651
$ws.Range("AQ222").Value2 = 13.3094401
$ws.Range("AS222").Value2 = 10.865627
$ws.Range("AT222").Value2 = 21.5665838
$ws.Range("AU222").Value2 = 33.2418035
$ws.Range("AV222").Value2 = 26.1858926
$ws.Range("AW222").Value2 = 21.2629632
$ws.Range("AX222").Value2 = 23.061841
$ws.Range("AY222").Value2 = 17.4436448
$ws.Range("BA222").Value2 = 7.7063163
$ws.Range("BB222").Value2 = 14.7090186
$ws.Range("BC222").Value2 = 19.6888794
$ws.Range("BD222").Value2 = 21.1989539
$ws.Range("BE222").Value2 = 21.0936722

# Row 223
$ws.Range("A223").Value2 = "09 09 2020"
$ws.Range("B223").Value2 = 19.8795181
$ws.Range("C223").Value2 = 26.5694584
$ws.Range("D223").Value2 = 26.5144293
$ws.Range("F223").Value2 = 16.5157421
$ws.Range("G223").Value2 = 14.5187308
$ws.Range("H223").Value2 = 15.9367374
$ws.Range("I223").Value2 = 9.3437006
$ws.Range("J223").Value2 = 12.300885
$ws.Range("K223").Value2 = 11.9318182
$ws.Range("L223").Value2 = 19.1864101
$ws.Range("M223").Value2 = 25.2825168
$ws.Range("O223").Value2 = 15.9129693
$ws.Range("P223").Value2 = 28.8272992
$ws.Range("Q223").Value2 = 27.6110537
$ws.Range("R223").Value2 = 20.3435622
$ws.Range("S223").Value2 = 24.1517566
$ws.Range("T223").Value2 = 25.4247589
$ws.Range("U223").Value2 = 22.8076706
$ws.Range("V223").Value2 = 23.0495445
$ws.Range("W223").Value2 = 10.339222
$ws.Range("X223").Value2 = 13.3773068
$ws.Range("Y223").Value2 = 11.1081081
$ws.Range("Z223").Value2 = 14.706838
$ws.Range("AA223").Value2 = 16.828343
$ws.Range("AB223").Value2 = 28.6941535
$ws.Range("AD223").Value2 = 28.0031041
$ws.Range("AE223").Value2 = 22.2854083
$ws.Range("AF223").Value2 = 17.5352234
$ws.Range("AG223").Value2 = 31.5016025
$ws.Range("AH223").Value2 = 24.990026
$ws.Range("AI223").Value2 = 8.5359116
$ws.Range("AJ223").Value2 = 10.0287742
$ws.Range("AK223").Value2 = 14.7268807
$ws.Range("AL223").Value2 = 18.8086993
$ws.Range("AM223").Value2 = 9.9657599
$ws.Range("AN223").Value2 = 17.4098234
$ws.Range("AO223").Value2 = 27.3153687
$ws.Range("AP223").Value2 = 12.5686644
$ws.Range("AQ223").Value2 = 13.1950802
$ws.Range("AS223").Value2 = 11.2301506
$ws.Range("AT223").Value2 = 21.407653
$ws.Range("AU223").Value2 = 33.1859061
$ws.Range("AV223").Value2 = 25.5625221
$ws.Range("AW223").Value2 = 20.8738051
$ws.Range("AX223").Value2 = 23.3509707
$ws.Range("AY223").Value2 = 17.7345454
$ws.Range("BA223").Value2 = 7.7651592
$ws.Range("BB223").Value2 = 14.7010333
$ws.Range("BC223").Value2 = 19.3623978
$ws.Range("BD223").Value2 = 21.0092047
$ws.Range("BE223").Value2 = 21.269023

# Row 224
$ws.Range("A224").Value2 = "10 09 2020"
$ws.Range("B224").Value2 = 20.0137552
$ws.Range("C224").Value2 = 26.2870897
$ws.Range("D224").Value2 = 26.5607652
$ws.Range("F224").Value2 = 16.6649902
$ws.Range("G224").Value2 = 14.4708911
$ws.Range("H224").Value2 = 16.0365701
$ws.Range("I224").Value2 = 9.4495662
$ws.Range("J224").Value2 = 10.6284658
$ws.Range("K224").Value2 = 10.7397504
$ws.Range("L224").Value2 = 19.1794418
$ws.Range("M224").Value2 = 25.2582446
$ws.Range("O224").Value2 = 16.5936953
$ws.Range("P224").Value2 = 28.1415347
$ws.Range("Q224").Value2 = 26.8073999
$ws.Range("R224").Value2 = 20.0992991
$ws.Range("S224").Value2 = 23.7292097
$ws.Range("T224").Value2 = 24.9662885
$ws.Range("U224").Value2 = 22.3454487
$ws.Range("V224").Value2 = 22.6882182
$ws.Range("W224").Value2 = 10.4502651
$ws.Range("X224").Value2 = 13.5671172
$ws.Range("Y224").Value2 = 10.7782755
$ws.Range("Z224").Value2 = 14.7280546
$ws.Range("AA224").Value2 = 17.1984238
$ws.Range("AB224").Value2 = 28.5267736
$ws.Range("AD224").Value2 = 27.1270979
$ws.Range("AE224").Value2 = 22.8381891
$ws.Range("AF224").Value2 = 17.39156
$ws.Range("AG224").Value2 = 31.9282043
$ws.Range("AH224").Value2 = 24.9791237
$ws.Range("AI224").Value2 = 8.581879600000001
$ws.Range("AJ224").Value2 = 9.786065300000001
$ws.Range("AK224").Value2 = 14.403464
$ws.Range("AL224").Value2 = 18.1072602
$ws.Range("AM224").Value2 = 9.847248
$ws.Range("AN224").Value2 = 16.8482128
$ws.Range("AO224").Value2 = 27.4766045
$ws.Range("AP224").Value2 = 12.3316526
$ws.Range("AQ224").Value2 = 13.3574794
$ws.Range("AS224").Value2 = 11.1482887
$ws.Range("AT224").Value2 = 22.2697392
$ws.Range("AU224").Value2 = 32.0536947
$ws.Range("AV224").Value2 = 25.8324932
$ws.Range("AW224").Value2 = 20.5434195
$ws.Range("AX224").Value2 = 23.4884273
$ws.Range("AY224").Value2 = 17.7084042
$ws.Range("BA224").Value2 = 7.0625099
$ws.Range("BB224").Value2 = 14.5224421
$ws.Range("BC224").Value2 = 19.4397107
$ws.Range("BD224").Value2 = 20.7332349
$ws.Range("BE224").Value2 = 20.2357396

# Row 225
$ws.Range("A225").Value2 = "11 09 2020"
$ws.Range("B225").Value2 = 20.6919946
$ws.Range("C225").Value2 = 26.3541035
$ws.Range("D225").Value2 = 26.4438674
$ws.Range("F225").Value2 = 16.6404644
$ws.Range("G225").Value2 = 14.1553496
$ws.Range("H225").Value2 = 15.7205043
$ws.Range("I225").Value2 = 9.4616863
$ws.Range("J225").Value2 = 9.7042514
$ws.Range("K225").Value2 = 10.8288066
$ws.Range("L225").Value2 = 19.0718649
$ws.Range("M225").Value2 = 25.0592481
$ws.Range("O225").Value2 = 17.5967597
$ws.Range("P225").Value2 = 27.6355026
$ws.Range("Q225").Value2 = 26.4948518
$ws.Range("R225").Value2 = 19.8771297
$ws.Range("S225").Value2 = 23.8275717
$ws.Range("T225").Value2 = 25.207784
$ws.Range("U225").Value2 = 22.071824
$ws.Range("V225").Value2 = 22.6937157
$ws.Range("W225").Value2 = 10.5130293
$ws.Range("X225").Value2 = 14.0519674
$ws.Range("Y225").Value2 = 11.1141907
$ws.Range("Z225").Value2 = 14.6446743
$ws.Range("AA225").Value2 = 17.6617248
$ws.Range("AB225").Value2 = 28.2102617
$ws.Range("AD225").Value2 = 26.658775
$ws.Range("AE225").Value2 = 23.4346398
$ws.Range("AF225").Value2 = 17.3000662
$ws.Range("AG225").Value2 = 31.0284026
$ws.Range("AH225").Value2 = 24.4604083
$ws.Range("AI225").Value2 = 8.4347826
$ws.Range("AJ225").Value2 = 9.8917913
$ws.Range("AK225").Value2 = 14.3413657
$ws.Range("AL225").Value2 = 17.7195963
$ws.Range("AM225").Value2 = 9.8332639
$ws.Range("AN225").Value2 = 16.8487054
$ws.Range("AO225").Value2 = 27.5088576
$ws.Range("AP225").Value2 = 11.9613438
$ws.Range("AQ225").Value2 = 13.1351518
$ws.Range("AS225").Value2 = 10.4043298
$ws.Range("AT225").Value2 = 22.119396
$ws.Range("AU225").Value2 = 32.1527431
$ws.Range("AV225").Value2 = 25.7001976
$ws.Range("AW225").Value2 = 20.4386306
$ws.Range("AX225").Value2 = 23.6783527
$ws.Range("AY225").Value2 = 17.848186
$ws.Range("BA225").Value2 = 6.7848402
$ws.Range("BB225").Value2 = 14.3813926
$ws.Range("BC225").Value2 = 19.3631759
$ws.Range("BD225").Value2 = 21.0591425
$ws.Range("BE225").Value2 = 19.1395063

# Row 226
$ws.Range("A226").Value2 = "12 09 2020"
$ws.Range("B226").Value2 = 21.4189189
$ws.Range("C226").Value2 = 25.8250591
$ws.Range("D226").Value2 = 25.5466055
$ws.Range("F226").Value2 = 16.2977445
$ws.Range("G226").Value2 = 14.0317662
$ws.Range("H226").Value2 = 15.4612885
$ws.Range("I226").Value2 = 9.5284327
$ws.Range("J226").Value2 = 10.2398524
$ws.Range("K226").Value2 = 10.8772079
$ws.Range("L226").Value2 = 19.0430499
$ws.Range("M226").Value2 = 25.086284
$ws.Range("O226").Value2 = 17.8240741
$ws.Range("P226").Value2 = 27.3059261
$ws.Range("Q226").Value2 = 26.1393656
$ws.Range("R226").Value2 = 19.5059466
$ws.Range("S226").Value2 = 23.7792627
$ws.Range("T226").Value2 = 24.7940781
$ws.Range("U226").Value2 = 22.1743276
$ws.Range("V226").Value2 = 22.2856268
$ws.Range("W226").Value2 = 10.6389635
$ws.Range("X226").Value2 = 13.619983
$ws.Range("Y226").Value2 = 11.0216546
$ws.Range("Z226").Value2 = 14.7444598
$ws.Range("AA226").Value2 = 17.4469844
$ws.Range("AB226").Value2 = 28.1306928
$ws.Range("AD226").Value2 = 26.2499178
$ws.Range("AE226").Value2 = 22.6673065
$ws.Range("AF226").Value2 = 17.7030127
$ws.Range("AG226").Value2 = 31.2354901
$ws.Range("AH226").Value2 = 23.749632
$ws.Range("AI226").Value2 = 8.228211
$ws.Range("AJ226").Value2 = 9.7955822
$ws.Range("AK226").Value2 = 14.1052155
$ws.Range("AL226").Value2 = 17.9967298
$ws.Range("AM226").Value2 = 9.806152300000001
$ws.Range("AN226").Value2 = 16.6483822
$ws.Range("AO226").Value2 = 27.0264588
$ws.Range("AP226").Value2 = 12.1285944
$ws.Range("AQ226").Value2 = 13.1523162
$ws.Range("AS226").Value2 = 10.0101272
$ws.Range("AT226").Value2 = 21.3167706
$ws.Range("AU226").Value2 = 31.2205508
$ws.Range("AV226").Value2 = 25.6877983
$ws.Range("AW226").Value2 = 20.1524895
$ws.Range("AX226").Value2 = 23.3966747
$ws.Range("AY226").Value2 = 17.8303427
$ws.Range("BA226").Value2 = 6.1186643
$ws.Range("BB226").Value2 = 14.3435682
$ws.Range("BC226").Value2 = 19.340073
$ws.Range("BD226").Value2 = 20.8662937
$ws.Range("BE226").Value2 = 19.4468629

# Row 227
$ws.Range("A227").Value2 = "13 09 2020"
$ws.Range("B227").Value2 = 21.1333333
$ws.Range("C227").Value2 = 25.4461216
$ws.Range("D227").Value2 = 25.1917575
$ws.Range("F227").Value2 = 15.698215
$ws.Range("G227").Value2 = 13.9543655
$ws.Range("H227").Value2 = 15.5845761
$ws.Range("I227").Value2 = 9.9583333
$ws.Range("J227").Value2 = 9.847036299999999
$ws.Range("K227").Value2 = 11.066523
$ws.Range("L227").Value2 = 18.7404372
$ws.Range("M227").Value2 = 24.6861913
$ws.Range("O227").Value2 = 17.8023033
$ws.Range("P227").Value2 = 26.4897553
$ws.Range("Q227").Value2 = 24.9924095
$ws.Range("R227").Value2 = 19.4919148
$ws.Range("S227").Value2 = 23.3217981
$ws.Range("T227").Value2 = 24.2969426
$ws.Range("U227").Value2 = 22.1003218
$ws.Range("V227").Value2 = 21.9120412
$ws.Range("W227").Value2 = 10.4720499
$ws.Range("X227").Value2 = 13.8073005
$ws.Range("Y227").Value2 = 11.1142534
$ws.Range("Z227").Value2 = 14.5962301
$ws.Range("AA227").Value2 = 17.3147612
$ws.Range("AB227").Value2 = 28.386176
$ws.Range("AD227").Value2 = 25.4164026
$ws.Range("AE227").Value2 = 22.3934198
$ws.Range("AF227").Value2 = 17.8503425
$ws.Range("AG227").Value2 = 31.7972962
$ws.Range("AH227").Value2 = 23.7028982
$ws.Range("AI227").Value2 = 8.0294118
$ws.Range("AJ227").Value2 = 9.607103499999999
$ws.Range("AK227").Value2 = 13.2481981
$ws.Range("AL227").Value2 = 16.8974704
$ws.Range("AM227").Value2 = 9.563927899999999
$ws.Range("AN227").Value2 = 16.8538602
$ws.Range("AO227").Value2 = 26.9852041
$ws.Range("AP227").Value2 = 12.3068733
$ws.Range("AQ227").Value2 = 13.0799619
$ws.Range("AS227").Value2 = 9.3976524
$ws.Range("AT227").Value2 = 21.2262332
$ws.Range("AU227").Value2 = 30.899325
$ws.Range("AV227").Value2 = 25.3589209
$ws.Range("AW227").Value2 = 19.6635631
$ws.Range("AX227").Value2 = 22.9577575
$ws.Range("AY227").Value2 = 17.7627881
$ws.Range("BA227").Value2 = 5.6662509
$ws.Range("BB227").Value2 = 14.3978399
$ws.Range("BC227").Value2 = 19.7371071
$ws.Range("BD227").Value2 = 21.0405897
$ws.Range("BE227").Value2 = 19.8245423

# Row 228
$ws.Range("A228").Value2 = "14 09 2020"
